$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1062.375
$ws.Range("I12").Value = 885.4286
$ws.Range("K12").Value = 885.4286
$ws.Range("M12").Value = -715.4286
$ws.Range("H17").Value = 649.092
$ws.Range("J17").Value = 633.3837
$ws.Range("L17").Value = 1900.1511
$ws.Range("N17").Value = -2236.1511
$ws.Range("H28").Value = 33333908
$ws.Range("I28").Value = 40000376
$ws.Range("J28").Value = 1574.2
$ws.Range("K28").Value = 40000376
$ws.Range("L28").Value = 1574.2
$ws.Range("M28").Value = -39999891
$ws.Range("N28").Value = -2544.2
$ws.Range("H41").Value = 290.5
$ws.Range("I41").Value = 205.6875
$ws.Range("J41").Value = 375.3125
$ws.Range("K41").Value = 205.6875
$ws.Range("L41").Value = 375.3125
$ws.Range("M41").Value = 234.3125
$ws.Range("N41").Value = -1255.3125
$ws.Range("H53").Value = 316.92856
$ws.Range("I53").Value = 456.66666
$ws.Range("K53").Value = 456.66666
$ws.Range("M53").Value = 180.33334
$ws.Range("H62").Value = 3157.9092
$ws.Range("I62").Value = 2592.25
$ws.Range("J62").Value = 4666.3335
$ws.Range("K62").Value = 2592.25
$ws.Range("L62").Value = 4666.3335
$ws.Range("M62").Value = -1968.25
$ws.Range("N62").Value = -5914.3335
$ws.Range("H65").Value = 3157.9092
$ws.Range("I65").Value = 2592.25
$ws.Range("J65").Value = 4666.3335
$ws.Range("K65").Value = 12961.25
$ws.Range("L65").Value = 23331.6675
$ws.Range("M65").Value = -9841.25
$ws.Range("N65").Value = -29571.6675
$ws.Range("H76").Value = 27868.334
$ws.Range("J76").Value = 34994
$ws.Range("L76").Value = 34994
$ws.Range("N76").Value = -35624
$ws.Range("H79").Value = 27868.334
$ws.Range("J79").Value = 34994
$ws.Range("L79").Value = 34994
$ws.Range("N79").Value = -37178
$ws.Range("H86").Value = 2203.923
$ws.Range("I86").Value = 2449.3
$ws.Range("J86").Value = 1386
$ws.Range("K86").Value = 2449.3
$ws.Range("L86").Value = 1386
$ws.Range("M86").Value = -1326.3
$ws.Range("N86").Value = -3632
$ws.Range("H89").Value = 2203.923
$ws.Range("I89").Value = 2449.3
$ws.Range("J89").Value = 1386
$ws.Range("K89").Value = 12246.5
$ws.Range("L89").Value = 6930
$ws.Range("M89").Value = -6630.5
$ws.Range("N89").Value = -18162
$ws.Range("H92").Value = 88.9375
$ws.Range("J92").Value = 105.5
$ws.Range("L92").Value = 105.5
$ws.Range("N92").Value = -2601.5
$ws.Range("H98").Value = 1438.7097
$ws.Range("I98").Value = 1413.6957
$ws.Range("J98").Value = 1510.625
$ws.Range("K98").Value = 1413.6957
$ws.Range("L98").Value = 1510.625
$ws.Range("M98").Value = 84.30430000000001
$ws.Range("N98").Value = -4506.625
$ws.Range("H106").Value = 83339730
$ws.Range("I106").Value = 111116296
$ws.Range("J106").Value = 10000
$ws.Range("K106").Value = 111116296
$ws.Range("L106").Value = 10000
$ws.Range("M106").Value = -111115665
$ws.Range("N106").Value = -11262
$ws.Range("H107").Value = 15627362
$ws.Range("I107").Value = 17243014
$ws.Range("J107").Value = 9399.333000000001
$ws.Range("K107").Value = 17243014
$ws.Range("L107").Value = 9399.333000000001
$ws.Range("M107").Value = -17241094
$ws.Range("N107").Value = -13239.333
$ws.Range("H113").Value = 5378.5713
$ws.Range("H122").Value = 1438.7097
$ws.Range("I122").Value = 1413.6957
$ws.Range("J122").Value = 1510.625
$ws.Range("K122").Value = 4241.0871
$ws.Range("L122").Value = 4531.875
$ws.Range("M122").Value = -1791.0871
$ws.Range("N122").Value = -9431.875
$ws.Range("H137").Value = 11492.692
$ws.Range("J137").Value = 2442.8333
$ws.Range("L137").Value = 7328.499899999999
$ws.Range("N137").Value = -12428.4999
$ws.Range("H140").Value = 108685.14
$ws.Range("J140").Value = 148975
$ws.Range("L140").Value = 148975
$ws.Range("N140").Value = -159335
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14683949
$ws.Range("I2").Value = 24027216
$ws.Range("K2").Value = 24027216
$ws.Range("M2").Value = -24027103
$ws.Range("H32").Value = 5031.0864
$ws.Range("I32").Value = 4839.551
$ws.Range("J32").Value = 6073.8887
$ws.Range("K32").Value = 4839.551
$ws.Range("L32").Value = 6073.8887
$ws.Range("M32").Value = -4552.551
$ws.Range("N32").Value = -6647.8887
$ws.Range("H46").Value = 9393
$ws.Range("J46").Value = 10750.2
$ws.Range("L46").Value = 10750.2
$ws.Range("N46").Value = -11388.2
$ws.Range("H61").Value = 3980.1025
$ws.Range("I61").Value = 3958.7778
$ws.Range("K61").Value = 3958.7778
$ws.Range("M61").Value = -3746.7778
$ws.Range("H74").Value = 7081.724
$ws.Range("I74").Value = 7081.724
$ws.Range("K74").Value = 7081.724
$ws.Range("M74").Value = -6207.724
$ws.Range("H77").Value = 7081.724
$ws.Range("I77").Value = 7081.724
$ws.Range("K77").Value = 35408.62
$ws.Range("M77").Value = -31040.62
$ws.Range("H116").Value = 14683949
$ws.Range("I116").Value = 24027216
$ws.Range("K116").Value = 24027216
$ws.Range("M116").Value = -24024922
$ws.Range("H136").Value = 3980.1025
$ws.Range("I136").Value = 3958.7778
$ws.Range("K136").Value = 11876.3334
$ws.Range("M136").Value = -9326.3334
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14683949
$ws.Range("I3").Value = 24027216
$ws.Range("K3").Value = 24027216
$ws.Range("M3").Value = -24027102
$ws.Range("H134").Value = 2742.0173
$ws.Range("I134").Value = 2691.152
$ws.Range("J134").Value = 2937
$ws.Range("K134").Value = 8073.456
$ws.Range("L134").Value = 8811
$ws.Range("M134").Value = -5538.456
$ws.Range("N134").Value = -13881
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 22650
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H69").Value = 39750
$ws.Range("I69").Value = 39750
$ws.Range("K69").Value = 39750
$ws.Range("M69").Value = -39001
$ws.Range("H72").Value = 39750
$ws.Range("I72").Value = 39750
$ws.Range("K72").Value = 119250
$ws.Range("M72").Value = -115506
$ws.Range("H138").Value = 144890.67
$ws.Range("J138").Value = 144890.67
$ws.Range("L138").Value = 144890.67
$ws.Range("N138").Value = -155170.67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 490.79166
$ws.Range("I5").Value = 389.3158
$ws.Range("J5").Value = 876.4
$ws.Range("K5").Value = 1167.9474
$ws.Range("L5").Value = 2629.2
$ws.Range("M5").Value = -1055.9474
$ws.Range("N5").Value = -2853.2
$ws.Range("H135").Value = 490.79166
$ws.Range("I135").Value = 389.3158
$ws.Range("J135").Value = 876.4
$ws.Range("K135").Value = 3503.8422
$ws.Range("L135").Value = 7887.599999999999
$ws.Range("M135").Value = -968.8422
$ws.Range("N135").Value = -12957.6
$ws.Range("H136").Value = 12584
$ws.Range("J136").Value = 17500
$ws.Range("L136").Value = 52500
$ws.Range("N136").Value = -62700
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 10134748
$ws.Range("I7").Value = 179663.33
$ws.Range("K7").Value = 179663.33
$ws.Range("M7").Value = -179551.33
$ws.Range("H8").Value = 10134748
$ws.Range("I8").Value = 179663.33
$ws.Range("K8").Value = 179663.33
$ws.Range("M8").Value = -179524.33
$ws.Range("H14").Value = 7254004
$ws.Range("I14").Value = 20000000
$ws.Range("J14").Value = 4067504.8
$ws.Range("K14").Value = 20000000
$ws.Range("L14").Value = 4067504.8
$ws.Range("M14").Value = -19999832
$ws.Range("N14").Value = -4067840.8
$ws.Range("H49").Value = 20000
$ws.Range("J49").Value = 20000
$ws.Range("L49").Value = 20000
$ws.Range("N49").Value = -20368
$ws.Range("H122").Value = 4051.2646
$ws.Range("I122").Value = 3646.6428
$ws.Range("K122").Value = 10939.9284
$ws.Range("M122").Value = -8489.928400000001
$ws.Range("H140").Value = 58965.2
$ws.Range("J140").Value = 69956.5
$ws.Range("L140").Value = 69956.5
$ws.Range("N140").Value = -80316.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4343.485
$ws.Range("I40").Value = 4537.773
$ws.Range("J40").Value = 3954.9092
$ws.Range("K40").Value = 4537.773
$ws.Range("L40").Value = 3954.9092
$ws.Range("M40").Value = -4401.773
$ws.Range("N40").Value = -4226.9092
$ws.Range("H46").Value = 4054.8125
$ws.Range("I46").Value = 2697.6
$ws.Range("J46").Value = 4671.727
$ws.Range("K46").Value = 2697.6
$ws.Range("L46").Value = 4671.727
$ws.Range("M46").Value = -2509.6
$ws.Range("N46").Value = -5047.727
$ws.Range("H92").Value = 67221.5
$ws.Range("J92").Value = 67221.5
$ws.Range("L92").Value = 67221.5
$ws.Range("N92").Value = -72213.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2757687.2
$ws.Range("I81").Value = 2935344.5
$ws.Range("J81").Value = 3998.5
$ws.Range("K81").Value = 5870689
$ws.Range("L81").Value = 7997
$ws.Range("M81").Value = -5869628
$ws.Range("N81").Value = -10119
$ws.Range("H84").Value = 2757687.2
$ws.Range("I84").Value = 2935344.5
$ws.Range("J84").Value = 3998.5
$ws.Range("K84").Value = 29353445
$ws.Range("L84").Value = 39985
$ws.Range("M84").Value = -29348141
$ws.Range("N84").Value = -50593
